$wb = $excel.ActiveWorkbook

# The workbook has 4 sheets: 展览 (sheet1), 演出 (sheet2), 本地生活 (sheet3), 全部类型 (sheet4).
# Rows for a canceled event appear identically in both "展览" and "全部类型" sheets (row 3),
# along with several "想去人数" (interest count, column F) bumps on other rows.

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 3: event cancelled - append marker to name, mark ticket price column as unavailable
    $ws.Range("C3").Value = "合肥·IE动漫嘉年华（取消）"
    $ws.Range("G3").Value = "不可售"

    # Update "interest" counts (column F) for several rows
    $ws.Range("F6").Value = 185
    $ws.Range("F10").Value = 5557
    $ws.Range("F11").Value = 4933
    $ws.Range("F13").Value = 43
}
